$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'" + '27.553.80'
$ws.Range("E2").Value = '  -0.23%  '

$ws.Range("D3").Value = "'" + '1.647.05'
$ws.Range("E3").Value = '  -0.65%  '

$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").Value = "'" + '212.54'
$ws.Range("E5").Value = '  -1.16%  '

$ws.Range("E6").Value = '  +4.99%  '

$ws.Range("E7").Value = '  +0.05%  '

$ws.Range("D8").Value = "'" + '23.57'
$ws.Range("E8").Value = '  -2.37%  '

$ws.Range("E9").Value = '  -1.93%  '

$ws.Range("E10").Value = '  -1.39%  '

$ws.Range("D11").Value = "'" + '0.0889'
$ws.Range("E11").Value = '  +1.01%  '

$ws.Range("D12").Value = "'" + '1.880.04'
$ws.Range("E12").Value = '  -0.72%  '

$ws.Range("D13").Value = "'" + '1.660.88'
$ws.Range("E13").Value = '  +0.32%  '

$ws.Range("D14").Value = "'" + '0.585'
$ws.Range("E14").Value = '  +3.12%  '

$ws.Range("E15").Value = '  -2.58%  '

$ws.Range("D16").Value = "'" + '64.50'
$ws.Range("E16").Value = '  -2.11%  '

$ws.Range("D17").Value = "'" + '27.526.12'
$ws.Range("E17").Value = '  -0.29%  '

$ws.Range("D18").Value = "'" + '230.43'
$ws.Range("E18").Value = '  -4.50%  '

$ws.Range("D19").Value = "'" + '0.0₃0724'
$ws.Range("E19").Value = '  -0.70%  '

$ws.Range("D20").Value = "'" + '7.55'
$ws.Range("E20").Value = '  -0.06%  '

$ws.Range("E21").Value = '  +0.12%  '

$ws.Range("E22").Value = '  -3.65%  '

$ws.Range("D23").Value = "'" + '9.74'
$ws.Range("E23").Value = '  +3.82%  '

$ws.Range("E24").Value = '  -2.28%  '

$ws.Range("D25").Value = "'" + '148.94'
$ws.Range("E25").Value = '  +2.04%  '

$ws.Range("D26").Value = "'" + '7.01'
$ws.Range("E26").Value = '  -2.82%  '

$ws.Range("E27").Value = '  +1.30%  '

$ws.Range("E28").Value = '  +0.00%  '

$ws.Range("D29").Value = "'" + '15.58'
$ws.Range("E29").Value = '  -4.49%  '

$ws.Range("E30").Value = '  -2.45%  '

$ws.Range("D31").Value = "'" + '0.0486'
$ws.Range("E31").Value = '  -3.29%  '

$ws.Range("E32").Value = '  -0.82%  '

$ws.Range("D33").Value = "'" + '3.19'
$ws.Range("E33").Value = '  +2.22%  '

$ws.Range("D34").Value = "'" + '1.425.59'
$ws.Range("E34").Value = '  -2.10%  '

$ws.Range("E35").Value = '  +1.37%  '

$ws.Range("D36").Value = "'" + '2.38'
$ws.Range("E36").Value = '  +0.05%  '

$ws.Range("D37").Value = "'" + '0.569'
$ws.Range("E37").Value = '  -0.92%  '

$ws.Range("E38").Value = '  -4.27%  '

$ws.Range("E39").Value = '  -2.98%  '

$ws.Range("E40").Value = '  +0.20%  '

$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D42").Value = "'" + '5.55'
$ws.Range("E42").Value = '  +2.27%  '

$ws.Range("B43").Value = 'TrustWalletToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D43").Value = "'" + '0.817'
$ws.Range("E43").Value = '  +2.90%  '

$ws.Range("E44").Value = '  +1.28%  '

$ws.Range("D45").Value = "'" + '65.09'
$ws.Range("E45").Value = '  -6.67%  '

$ws.Range("D46").Value = "'" + '1.789.63'
$ws.Range("E46").Value = '  -0.67%  '

$ws.Range("D47").Value = "'" + '1.68'
$ws.Range("E47").Value = '  -1.75%  '

$ws.Range("E48").Value = '  -0.54%  '

$ws.Range("E49").Value = '  +1.39%  '

$ws.Range("D50").Value = "'" + '0.0994'
$ws.Range("E50").Value = '  -3.30%  '

$ws.Range("D51").Value = "'" + '7.78'
$ws.Range("E51").Value = '  -0.70%  '
